{"js": "// The document's date line reads \"\u6295\u7968\u65f6\u95f4\uff1a2022\u5e743\u670827\u65e5\" (voting date).\n// The edit changes the day from the 27th to the 28th, i.e. the single\n// character run containing \"7\" (the day's units digit) becomes \"8\".\nconst body = context.document.body;\n\n// Scope the search to the paragraph that holds the voting-date text so the\n// replacement stays unambiguous even if the surrounding document changes.\nconst dateParagraphs = body.search(\"\u6295\u7968\u65f6\u95f4\uff1a*\u65e5\", { matchWildcards: true });\ndateParagraphs.load(\"items\");\nawait context.sync();\n\nconst searchScope = dateParagraphs.items.length > 0 ? dateParagraphs.items[0] : body;\n\nconst matches = searchScope.search(\"7\", { matchCase: true });\nmatches.load(\"items,text\");\nawait context.sync();\n\nif (matches.items.length !== 1) {\n  throw new Error(\"Expected exactly one '7' run in the voting-date text, found \" + matches.items.length);\n}\n\n// Replace just that run's text in place, preserving its run formatting\n// (font, size, rsid, etc.) exactly like editing the <w:t> content in OOXML.\nmatches.items[0].insertText(\"8\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document's voting-date line reads \"\u6295\u7968\u65f6\u95f4\uff1a2022\u5e743\u670827\u65e5\".\n# This edit changes the day from the 27th to the 28th, i.e. the character\n# \"7\" (the day's units digit) becomes \"8\".\n$d = $word.ActiveDocument\n\n# Sanity-check: make sure there is exactly one \"7\" in the document before\n# touching anything, so the replacement below is unambiguous.\n$probeRange = $d.Content.Duplicate\n$probeFind = $probeRange.Find\n$probeFind.ClearFormatting()\n$probeFind.Text = \"7\"\n$probeFind.MatchCase = $true\n$probeFind.Forward = $true\n$probeFind.Wrap = 0\n$matchCount = 0\nwhile ($probeFind.Execute()) {\n    $matchCount = $matchCount + 1\n}\nif ($matchCount -ne 1) {\n    throw \"Expected exactly one '7' to replace, found $matchCount.\"\n}\n\n# Standard Find & Replace (replace the first/only match) on the real content.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"7\"\n$find.Replacement.Text = \"8\"\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\n\n$found = $find.Execute($null, $true, $null, $null, $null, $null, $true, $null, $null, \"8\", 2)\nif (-not $found) {\n    throw \"Could not find the '7' in the voting-date text to replace.\"\n}\n"}
